$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold + border) from K1 onto the new L1 header cell,
# then set its text.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").Value = "Sizes"

# Sizes values for rows 2-70 (product rows), in row order.
$sizes = @(
    "N/A",
    "N/A",
    "N/A",
    "N/A",
    "N/A",
    "N/A",
    "129,139,149",
    "54MM",
    "8 1/2",
    "8 1/8",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge",
    "30,32,34,36,38",
    "30,32,34,36,38",
    "30,32,34,36,38",
    "Small,Medium,Large,XLarge,XXL",
    "30,32,34,36,38",
    "Small,Medium,Large,XLarge,XXL",
    "30,32,34,36,38",
    "30,32,34,36,38",
    "30,32,34,36,38",
    "30,32,34,36,38",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "N/A",
    "N/A",
    "N/A",
    "S/M,M/L",
    "N/A",
    "N/A",
    "N/A",
    "N/A",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "Small,Medium,Large,XLarge,XXL",
    "N/A"
)

for ($i = 0; $i -lt $sizes.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $sizes[$i]
}
